$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix comma-separated names to period-separated (OCR/scrape artifact fix)
$ws.Range("E94").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E125").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E220").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E205").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E215").Value = "OLIVERA. FLORENCIO"
$ws.Range("F215").Value = "OLIVERA. FLORENCIO"
$ws.Range("E217").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E238").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E264").Value = "ALVAREZ. RENZO JOEL"
$ws.Range("F264").Value = "ALVAREZ. RENZO JOEL"
$ws.Range("E268").Value = "ODIARD. OSCAR HERNAN"
$ws.Range("F268").Value = "ODIARD. OSCAR HERNAN"
$ws.Range("F129").Value = "PARRAVICINI VIRGINIA VANINA. VIRGINIA VANINA"
$ws.Range("F233").Value = "PARRAVICINI VIRGINIA VANINA. VIRGINIA VANINA"

# Fix floating point Importe values: "1.234.567,89" (Spanish format, stored as text)
# -> "1234567.89" (plain decimal dot, stored as text) - apostrophe forces text
# so Excel does not silently coerce the numeric-looking string into a Number cell.
$importeMap = @{
    2 = "709000.00"
    3 = "20000.00"
    4 = "960000.00"
    5 = "4450.00"
    6 = "990.00"
    7 = "53475.00"
    8 = "8840.00"
    9 = "44226.00"
    10 = "2844.07"
    11 = "1174200.60"
    12 = "3694.94"
    13 = "14005.29"
    14 = "2726.40"
    15 = "70070.00"
    16 = "22400.00"
    17 = "390730.45"
    18 = "4800.00"
    19 = "1871814.36"
    20 = "65560.00"
    21 = "322520.00"
    22 = "402811.62"
    23 = "8400.00"
    24 = "193327.94"
    25 = "115822.72"
    26 = "3529.40"
    27 = "58382.28"
    28 = "5995.00"
    29 = "66300.64"
    30 = "8800.00"
    31 = "5000.00"
    32 = "21780.00"
    33 = "500.94"
    34 = "290.00"
    35 = "137387.84"
    36 = "101264.00"
    37 = "817550.00"
    38 = "32070.00"
    39 = "26939.80"
    40 = "2100.00"
    41 = "9093.00"
    42 = "1571.76"
    43 = "658600.00"
    44 = "5227.50"
    45 = "211900.00"
    46 = "1518.00"
    47 = "85872.85"
    48 = "59000.00"
    49 = "19509.12"
    50 = "38244.00"
    51 = "1590.00"
    52 = "3400.00"
    53 = "103196.50"
    54 = "12568.97"
    55 = "1237.00"
    56 = "25982.36"
    57 = "6800.00"
    58 = "21411.51"
    59 = "21283.02"
    60 = "26000.00"
    61 = "590625.20"
    62 = "2530.00"
    63 = "2324.30"
    64 = "22799.00"
    65 = "1080.00"
    66 = "67056.45"
    67 = "20700.00"
    68 = "9456.64"
    69 = "11500.00"
    70 = "1257.00"
    71 = "722988.00"
    72 = "13600.00"
    73 = "8680.00"
    74 = "76430.50"
    75 = "2720.00"
    76 = "9810.00"
    77 = "22030.51"
    78 = "9810.00"
    79 = "38130.00"
    80 = "43480.00"
    81 = "53849.98"
    82 = "700.00"
    83 = "250000.00"
    84 = "1790.00"
    85 = "66600.00"
    86 = "21100.00"
    87 = "89202.00"
    88 = "1085.70"
    89 = "980.00"
    90 = "3300.00"
    91 = "3078.00"
    92 = "137775.00"
    93 = "46608.00"
    94 = "600.00"
    95 = "10450.00"
    96 = "22000.00"
    97 = "2200.00"
    98 = "36.64"
    99 = "10720.00"
    100 = "17106.35"
    101 = "18600.00"
    102 = "1000.00"
    103 = "3855.00"
    104 = "18965.02"
    105 = "14204.61"
    106 = "36173.43"
    107 = "10000.00"
    108 = "40091.00"
    109 = "15170.00"
    110 = "6580.00"
    111 = "22645.00"
    112 = "13305.00"
    113 = "759.80"
    114 = "1600.00"
    115 = "2510.00"
    116 = "17520.00"
    117 = "6850.00"
    118 = "5635.00"
    119 = "71685.68"
    120 = "2872.00"
    121 = "640.00"
    122 = "7320.00"
    123 = "1507.80"
    124 = "24957.00"
    125 = "1450.00"
    126 = "6200.00"
    127 = "44100.00"
    128 = "12250.00"
    129 = "5900.00"
    130 = "9323.00"
    131 = "7008.68"
    132 = "23700.00"
    133 = "5400.00"
    134 = "600.00"
    135 = "218991.00"
    136 = "8793.39"
    137 = "95760.00"
    138 = "5987.18"
    139 = "86000.00"
    140 = "11280.00"
    141 = "18000.00"
    142 = "4800.00"
    143 = "8000.00"
    144 = "5000.00"
    145 = "10000.00"
    146 = "9000.00"
    147 = "5000.00"
    148 = "13849.00"
    149 = "22000.00"
    150 = "41000.00"
    151 = "2000.00"
    152 = "22000.00"
    153 = "78000.00"
    154 = "19000.00"
    155 = "39000.00"
    156 = "7000.20"
    157 = "17000.00"
    158 = "12637.47"
    159 = "32974.92"
    160 = "19376.00"
    161 = "8692.77"
    162 = "80181.00"
    163 = "28840.00"
    164 = "28360.04"
    165 = "10719.24"
    166 = "68110.00"
    167 = "98620.00"
    168 = "1643.40"
    169 = "2940000.00"
    170 = "6050.00"
    171 = "1080.00"
    172 = "42000.00"
    173 = "16500.00"
    174 = "22000.00"
    175 = "22000.00"
    176 = "20000.00"
    177 = "7000.00"
    178 = "22000.00"
    179 = "52000.00"
    180 = "10000.00"
    181 = "12000.00"
    182 = "10000.00"
    183 = "10000.00"
    184 = "9000.00"
    185 = "10500.00"
    186 = "8000.00"
    187 = "10000.00"
    188 = "6000.00"
    189 = "10000.00"
    190 = "22000.00"
    191 = "10000.00"
    192 = "10000.00"
    193 = "20870.00"
    194 = "5000.00"
    195 = "45000.00"
    196 = "18000.00"
    197 = "10000.00"
    198 = "12000.00"
    199 = "25000.00"
    200 = "7000.00"
    201 = "5000.00"
    202 = "23205.00"
    203 = "8000.00"
    204 = "20500.00"
    205 = "10000.00"
    206 = "5000.00"
    207 = "10000.00"
    208 = "10000.00"
    209 = "58000.00"
    210 = "36000.00"
    211 = "45000.00"
    212 = "9800.00"
    213 = "29300.00"
    214 = "19870.00"
    215 = "5000.00"
    216 = "46700.00"
    217 = "2800.00"
    218 = "2200.00"
    219 = "7974.15"
    220 = "20175.00"
    221 = "39750.00"
    222 = "2143.00"
    223 = "53120.00"
    224 = "43936.00"
    225 = "6900.00"
    226 = "4747.28"
    227 = "220.30"
    228 = "7180.00"
    229 = "150340.00"
    230 = "13200.00"
    231 = "11400.00"
    232 = "18264.84"
    233 = "19744.00"
    234 = "24688.82"
    235 = "30000.00"
    236 = "3185.00"
    237 = "1806.11"
    238 = "104130.00"
    239 = "9907.82"
    240 = "6100.00"
    241 = "7229.12"
    242 = "3200.00"
    243 = "16369.50"
    244 = "6408.60"
    245 = "4100.00"
    246 = "7500.00"
    247 = "19200.00"
    248 = "101300.00"
    249 = "15000.00"
    250 = "300.00"
    251 = "78907.47"
    252 = "1176.25"
    253 = "80000.00"
    254 = "40000.00"
    255 = "40000.00"
    256 = "40000.00"
    257 = "80000.00"
    258 = "40000.00"
    259 = "55000.00"
    260 = "40000.00"
    261 = "40000.00"
    262 = "80000.00"
    263 = "80000.00"
    264 = "23500.00"
    265 = "8500.00"
    266 = "14133.37"
    267 = "23700.00"
    268 = "14330.00"
    269 = "177600.00"
    270 = "14630.00"
    271 = "446000.00"
    272 = "59000.00"
    273 = "2534974.00"
    274 = "408210.00"
    275 = "372480.00"
    276 = "325800.00"
    277 = "347000.00"
    278 = "309000.00"
    279 = "309000.00"
    280 = "588000.00"
    281 = "309000.00"
    282 = "837990.00"
    283 = "1046000.00"
    284 = "365840.00"
    285 = "309000.00"
    286 = "309000.00"
    287 = "618000.00"
    288 = "638340.00"
    289 = "609840.00"
    290 = "910680.00"
    291 = "588000.00"
    292 = "898060.00"
    293 = "618000.00"
    294 = "337043.00"
    295 = "343239.76"
    296 = "956305.81"
    297 = "22400.00"
    298 = "3157649.87"
    299 = "3365889.91"
    300 = "1336.96"
    301 = "148780.95"
    302 = "8000.00"
    303 = "1960.00"
    304 = "4700.00"
    305 = "245340.00"
    306 = "45000.00"
    307 = "25192.17"
    308 = "1700.00"
    309 = "2929.52"
    310 = "163312.00"
}
foreach ($row in $importeMap.Keys) {
    $ws.Cells.Item($row, 8).Value = "'" + $importeMap[$row]
}
